$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(8)

# Start clean: clear all existing content/formatting on the sheet so we can rebuild it
$ws.Cells.Clear()

# --- Header row (row 1) ---
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# --- Data rows (rows 2-6) ---
# Row 2
$ws.Range("A2").Value = 140
$ws.Range("B2").Value = "三商美邦人壽"
$ws.Range("C2").Value = "祥安終身壽險"
$ws.Range("D2").Value = "黃停哲"
$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"
$ws.Range("H2").Value = "黃偉哲"
$ws.Range("J2").Value = "tmp85f1"
$ws.Range("G2").Value = "'2011-12-20"
$ws.Range("I2").Value = 1367
$ws.Range("K2").Value = 140

# Row 3
$ws.Range("A3").Value = 141
$ws.Range("B3").Value = "國泰人壽"
$ws.Range("C3").Value = "美滿人生202終身"
$ws.Range("D3").Value = "黃偉哲"
$ws.Range("E3").Value = "insurance"
$ws.Range("F3").Value = "normal"
$ws.Range("H3").Value = "黃偉哲"
$ws.Range("J3").Value = "tmp85f1"
$ws.Range("G3").Value = "'2011-12-20"
$ws.Range("I3").Value = 1367
$ws.Range("K3").Value = 141

# Row 4
$ws.Range("A4").Value = 142
$ws.Range("B4").Value = "國泰人壽"
$ws.Range("C4").Value = "住院醫療健康保險"
$ws.Range("D4").Value = "黃偉哲"
$ws.Range("E4").Value = "insurance"
$ws.Range("F4").Value = "normal"
$ws.Range("H4").Value = "黃偉哲"
$ws.Range("J4").Value = "tmp85f1"
$ws.Range("G4").Value = "'2011-12-20"
$ws.Range("I4").Value = 1367
$ws.Range("K4").Value = 142

# Row 5
$ws.Range("A5").Value = 143
$ws.Range("B5").Value = "國泰人壽"
$ws.Range("C5").Value = "美意年年終身年金型保單"
$ws.Range("D5").Value = "劉育菁"
$ws.Range("E5").Value = "insurance"
$ws.Range("F5").Value = "normal"
$ws.Range("H5").Value = "黃偉哲"
$ws.Range("J5").Value = "tmp85f1"
$ws.Range("G5").Value = "'2011-12-20"
$ws.Range("I5").Value = 1367
$ws.Range("K5").Value = 143

# Row 6
$ws.Range("A6").Value = 144
$ws.Range("B6").Value = "國泰人壽"
$ws.Range("C6").Value = "住院醫療健康保險"
$ws.Range("D6").Value = "黃〇文"
$ws.Range("E6").Value = "insurance"
$ws.Range("F6").Value = "normal"
$ws.Range("H6").Value = "黃偉哲"
$ws.Range("J6").Value = "tmp85f1"
$ws.Range("G6").Value = "'2011-12-20"
$ws.Range("I6").Value = 1367
$ws.Range("K6").Value = 144

# --- Styling ---
# Header row (B1:K1): bold, thin border, centered/top-aligned -- matches the header style used
# throughout the rest of this workbook (e.g. sheet "具有相當價值之財產").
$hdr = $ws.Range("B1:K1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Column A (the numeric index column) uses the same bold/bordered/centered style as the header,
# matching the convention used in every other sheet of this workbook. A1 itself stays empty.
$colA = $ws.Range("A2:A6")
$colA.Font.Bold = $true
$colA.Borders.LineStyle = 1
$colA.HorizontalAlignment = -4108
$colA.VerticalAlignment = -4160

$excel.CutCopyMode = $false
